# StockForm update: replace the sample stock-card import row with the
# REMA iPhone 11 Açık Mavi product, add new maliyet/maliyetKur/vatRate
# columns, reorder warehouse/quantity columns earlier and re-sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (row 1) - columns A..O keep their original meaning, so
#    only their position-dependent neighbours (P,Q onward) need to move.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "productCode"
$ws.Range("B1").Value = "productName"
$ws.Range("C1").Value = "unit"
$ws.Range("D1").Value = "shortDescription"
$ws.Range("E1").Value = "description"
$ws.Range("F1").Value = "companyCode"
$ws.Range("G1").Value = "branchCode"
$ws.Range("H1").Value = "gtip"
$ws.Range("I1").Value = "pluCode"
$ws.Range("J1").Value = "desi"
$ws.Range("K1").Value = "adetBoleni"
$ws.Range("L1").Value = "siraNo"
$ws.Range("M1").Value = "raf"
$ws.Range("N1").Value = "karMarji"
$ws.Range("O1").Value = "riskQuantities"
$ws.Range("P1").Value = "maliyet"
$ws.Range("Q1").Value = "maliyetKur"
$ws.Range("R1").Value = "warehouseName"
$ws.Range("S1").Value = "quantity"
$ws.Range("T1").Value = "brandName"
$ws.Range("U1").Value = "productType"
$ws.Range("V1").Value = "categories"
$ws.Range("W1").Value = "attributes"
$ws.Range("X1").Value = "vatRate"
$ws.Range("Y1").Value = "prices"
$ws.Range("Z1").Value = "barcodes"
$ws.Range("AA1").Value = "manufacturerName"
$ws.Range("AB1").Value = "marketNames"

# New header cells (Z1:AB1) need the same bold/border/centered style as
# the rest of row 1 - copy it from an already-styled header cell so the
# existing style index is reused instead of creating new ones.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("Z1:AB1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Data row (row 2) - new product values.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "REMA/iP-11/AçıkMavi"
$ws.Range("B2").Value = "REMA iPhone 11 Açık Mavi"
$ws.Range("C2").Value = "Adet"
$ws.Range("D2").Value = "Apple Iphone 11 Magsafe Wireless Şarj Özellikli Silikon 2mm Kamera Çıkıntılı Rema Kılıf"
$ws.Range("E2").Value = "Apple Iphone 11 Magsafe Wireless Şarj Özellikli Silikon 2mm Kamera Çıkıntılı Rema Kılıf"
$ws.Range("F2").Value = "VIP"
$ws.Range("G2").Value = "ETC"
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1
$ws.Range("N2").Value = 20
$ws.Range("O2").Value = 50
$ws.Range("P2").Value = 0.8
$ws.Range("P2").NumberFormat = "0.00"
$ws.Range("Q2").Value = "USD"
$ws.Range("R2").Value = "E-Ticaret"
$ws.Range("S2").Value = 50
$ws.Range("T2").Value = "VipCase"
$ws.Range("U2").Value = "BasitUrun"
$ws.Range("V2").Value = "Rema Kılıf"
$ws.Range("W2").Value = "Renk=Açık Mavi,Cihaz Modeli=iPhone 11"
$ws.Range("X2").Value = 20
$ws.Range("Y2").Value = "Tip 1 Bayiler=1,Tip 2 Bayiler=1.4,Tip 3 Bayiler=1.8,Maliyet=0.8,Perakende=150"
$ws.Range("Z2").Value = "REMA/iP-11/AçıkMavi,8683606399265"
$ws.Range("AB2").Value = "EGE"

# ---------------------------------------------------------------------
# 3. Column widths (only columns whose width actually changes need to
#    be touched; the rest keep the width they already had).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.5
$ws.Columns.Item(2).ColumnWidth = 30
$ws.Columns.Item(4).ColumnWidth = 75
$ws.Columns.Item(5).ColumnWidth = 75
$ws.Columns.Item(16).ColumnWidth = 11.166666666666666
$ws.Columns.Item(17).ColumnWidth = 11.166666666666666
$ws.Columns.Item(18).ColumnWidth = 13.666666666666666
$ws.Columns.Item(19).ColumnWidth = 6.833333333333333
$ws.Columns.Item(20).ColumnWidth = 9.5
$ws.Columns.Item(21).ColumnWidth = 10.166666666666666
$ws.Columns.Item(22).ColumnWidth = 8.333333333333334
$ws.Columns.Item(23).ColumnWidth = 40.833333333333336
$ws.Columns.Item(24).ColumnWidth = 40.833333333333336
$ws.Columns.Item(25).ColumnWidth = 61
$ws.Columns.Item(26).ColumnWidth = 40.166666666666664
$ws.Columns.Item(27).ColumnWidth = 15.666666666666666
$ws.Columns.Item(28).ColumnWidth = 11.333333333333334

# ---------------------------------------------------------------------
# 4. View: zoom to 150% and move the selection.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("B7").Select()

# ---------------------------------------------------------------------
# 5. Sort the data range by productCode (adds sortState to the sheet).
# ---------------------------------------------------------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:AB2"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()
